# Auto-generated Excel COM-interop script
# Applies a scheduled market-data refresh to the per-job Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). For each affected row this updates
# the computed market/price/profit columns (H..N):
#   H = currentAveragePrice        I = currentAveragePriceNQ
#   J = currentAveragePriceHQ      K = LevePriceNQ
#   L = LevePriceHQ                M = LeveProfitNQ
#   N = LeveProfitHQ
# Some rows have no NQ or HQ profit figure at all (e.g. when there is no NQ/HQ
# recipe), so cells in M/N are cleared rather than set when the source data drops them.

$wb = $excel.ActiveWorkbook


# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 3602.2222
$ws.Range("J69").Value = 3602.2222
$ws.Range("L69").Value = 10806.6666
$ws.Range("N69").Value = -12554.6666

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 3602.2222
$ws.Range("J72").Value = 3602.2222
$ws.Range("L72").Value = 32419.9998
$ws.Range("N72").Value = -41155.99980000001

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 2867.3333
$ws.Range("I86").Value = 2867
$ws.Range("J86").Value = 2867.6667
$ws.Range("K86").Value = 2867
$ws.Range("L86").Value = 2867.6667
$ws.Range("M86").Value = -1744
$ws.Range("N86").Value = -5113.6667

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 2867.3333
$ws.Range("I89").Value = 2867
$ws.Range("J89").Value = 2867.6667
$ws.Range("K89").Value = 14335
$ws.Range("L89").Value = 14338.3335
$ws.Range("M89").Value = -8719
$ws.Range("N89").Value = -25570.3335

# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 251926.75
$ws.Range("I98").Value = 251926.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 251926.75
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -250428.75
$ws.Range("N98").ClearContents()

# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 251926.75
$ws.Range("I122").Value = 251926.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 755780.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -753330.25
$ws.Range("N122").ClearContents()

# Row 123: Nearly Bare / Gaja Grimoire
$ws.Range("H123").Value = 22990
$ws.Range("J123").Value = 22990
$ws.Range("L123").Value = 22990
$ws.Range("N123").Value = -32790

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 20036
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 20036
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 180324
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -185244

# Row 127: Liquid Competence / Competent Craftsman's Draught
$ws.Range("H127").Value = 1906.25
$ws.Range("I127").Value = 1625
$ws.Range("J127").Value = 2000
$ws.Range("K127").Value = 4875
$ws.Range("L127").Value = 6000
$ws.Range("M127").Value = 85
$ws.Range("N127").Value = -15920


# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 7078.0586
$ws.Range("I32").Value = 3072.8472
$ws.Range("J32").Value = 29260.77
$ws.Range("K32").Value = 3072.8472
$ws.Range("L32").Value = 29260.77
$ws.Range("M32").Value = -2785.8472
$ws.Range("N32").Value = -29834.77

# Row 82: Belle of the Brawl / Titanium Vambraces of Fending
$ws.Range("H82").Value = 10181
$ws.Range("J82").Value = 10181
$ws.Range("L82").Value = 10181
$ws.Range("N82").Value = -10903

# Row 85: Shouldering the Shut-ins (L) / Titanium Vambraces of Fending
$ws.Range("H85").Value = 10181
$ws.Range("J85").Value = 10181
$ws.Range("L85").Value = 10181
$ws.Range("N85").Value = -12677

# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 703.0833
$ws.Range("I97").Value = 703.0833
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 703.0833
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -207.0833
$ws.Range("N97").ClearContents()

# Row 117: Signed, Shield, Delivered / Titanbronze Tower Shield
$ws.Range("H117").Value = 39000
$ws.Range("J117").Value = 39000
$ws.Range("L117").Value = 39000
$ws.Range("N117").Value = -48178

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2523.7222
$ws.Range("I122").Value = 2245.3845
$ws.Range("J122").Value = 2681.0435
$ws.Range("K122").Value = 6736.1535
$ws.Range("L122").Value = 8043.130500000001
$ws.Range("M122").Value = -4286.1535
$ws.Range("N122").Value = -12943.1305

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1514.2593
$ws.Range("I132").Value = 1154.0857
$ws.Range("J132").Value = 2177.7368
$ws.Range("K132").Value = 3462.2571
$ws.Range("L132").Value = 6533.2104
$ws.Range("M132").Value = -932.2571000000003
$ws.Range("N132").Value = -11593.2104


# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1532.0667
$ws.Range("I86").Value = 1390.0834
$ws.Range("J86").Value = 2100
$ws.Range("K86").Value = 1390.0834
$ws.Range("L86").Value = 2100
$ws.Range("M86").Value = -267.0834
$ws.Range("N86").Value = -4346

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1532.0667
$ws.Range("I89").Value = 1390.0834
$ws.Range("J89").Value = 2100
$ws.Range("K89").Value = 6950.416999999999
$ws.Range("L89").Value = 10500
$ws.Range("M89").Value = -1334.416999999999
$ws.Range("N89").Value = -21732

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1107.8966
$ws.Range("I105").Value = 1030.4166
$ws.Range("J105").Value = 1479.8
$ws.Range("K105").Value = 1030.4166
$ws.Range("L105").Value = 1479.8
$ws.Range("M105").Value = 716.5834
$ws.Range("N105").Value = -4973.8

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1009.8939
$ws.Range("I134").Value = 725.2727
$ws.Range("J134").Value = 1579.1364
$ws.Range("K134").Value = 2175.8181
$ws.Range("L134").Value = 4737.4092
$ws.Range("M134").Value = 359.1819
$ws.Range("N134").Value = -9807.4092


# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2291.8723
$ws.Range("I31").Value = 1340.8823
$ws.Range("J31").Value = 4779.077
$ws.Range("K31").Value = 1340.8823
$ws.Range("L31").Value = 4779.077
$ws.Range("M31").Value = -1045.8823
$ws.Range("N31").Value = -5369.077

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2291.8723
$ws.Range("I34").Value = 1340.8823
$ws.Range("J34").Value = 4779.077
$ws.Range("K34").Value = 1340.8823
$ws.Range("L34").Value = 4779.077
$ws.Range("M34").Value = -1138.8823
$ws.Range("N34").Value = -5183.077

# Row 127: In Rod We Trust / Red Pine Fishing Rod
$ws.Range("H127").Value = 50000
$ws.Range("J127").Value = 50000
$ws.Range("L127").Value = 50000
$ws.Range("N127").Value = -59920

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1093.3462
$ws.Range("I132").Value = 773.1818
$ws.Range("J132").Value = 2854.25
$ws.Range("K132").Value = 2319.5454
$ws.Range("L132").Value = 8562.75
$ws.Range("M132").Value = 210.4546
$ws.Range("N132").Value = -13622.75


# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 126: Imperial Palate / Glory Be Soup
$ws.Range("H126").Value = 4221.6665
$ws.Range("I126").Value = 1015
$ws.Range("K126").Value = 3045
$ws.Range("M126").Value = 1895

# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 1037.3158
$ws.Range("I129").Value = 496.27274
$ws.Range("J129").Value = 1781.25
$ws.Range("K129").Value = 1488.81822
$ws.Range("L129").Value = 5343.75
$ws.Range("M129").Value = 3511.18178
$ws.Range("N129").Value = -15343.75


# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 2260.6
$ws.Range("I80").Value = 2283.3333
$ws.Range("J80").Value = 2226.5
$ws.Range("K80").Value = 2283.3333
$ws.Range("L80").Value = 2226.5
$ws.Range("M80").Value = -1285.3333
$ws.Range("N80").Value = -4222.5

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 2260.6
$ws.Range("I83").Value = 2283.3333
$ws.Range("J83").Value = 2226.5
$ws.Range("K83").Value = 11416.6665
$ws.Range("L83").Value = 11132.5
$ws.Range("M83").Value = -6424.666499999999
$ws.Range("N83").Value = -21116.5

# Row 87: Embroiling Embroidery / Griffin Talon Needle
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90: The Lovely Hands of Haillenarte (L) / Griffin Talon Needle
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 109: You're My Wonderhall / Hematite Earrings of Healing
$ws.Range("H109").Value = 13642.5
$ws.Range("J109").Value = 13642.5
$ws.Range("L109").Value = 13642.5
$ws.Range("N109").Value = -15722.5

# Row 119: Bulking Up / Dwarven Mythril Rapier
$ws.Range("H119").Value = 36086.832
$ws.Range("J119").Value = 36086.832
$ws.Range("L119").Value = 36086.832
$ws.Range("N119").Value = -45762.832


# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 36: Campaign in the Membrane / Toadskin Jacket
$ws.Range("H36").Value = 32900
$ws.Range("J36").Value = 32900
$ws.Range("L36").Value = 32900
$ws.Range("N36").Value = -34024

# Row 100: Tiger in the Sack / Tiger Leather
$ws.Range("H100").Value = 2375
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2750
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2750
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3832

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 1268.6724
$ws.Range("I132").Value = 952.0465
$ws.Range("J132").Value = 2176.3333
$ws.Range("K132").Value = 2856.1395
$ws.Range("L132").Value = 6528.999899999999
$ws.Range("M132").Value = -326.1395000000002
$ws.Range("N132").Value = -11588.9999


# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 6231.64
$ws.Range("I81").Value = 11063.363
$ws.Range("J81").Value = 2435.2856
$ws.Range("K81").Value = 22126.726
$ws.Range("L81").Value = 4870.5712
$ws.Range("M81").Value = -21065.726
$ws.Range("N81").Value = -6992.5712

# Row 84: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 6231.64
$ws.Range("I84").Value = 11063.363
$ws.Range("J84").Value = 2435.2856
$ws.Range("K84").Value = 110633.63
$ws.Range("L84").Value = 24352.856
$ws.Range("M84").Value = -105329.63
$ws.Range("N84").Value = -34960.856

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 752.84906
$ws.Range("I132").Value = 514.02563
$ws.Range("J132").Value = 1418.1428
$ws.Range("K132").Value = 1542.07689
$ws.Range("L132").Value = 4254.428400000001
$ws.Range("M132").Value = 987.9231100000002
$ws.Range("N132").Value = -9314.428400000001
